$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 7.050186157226562
$ws.Range("B1").Value = 5.941616058349609
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.244271755218506
$ws.Range("E1").Value = 1.910949468612671
